{"js": "// Insert a new paragraph \"Commit pinkesh 2\" right after the existing\n// \"Commit pinkesh\" paragraph (and before the trailing empty paragraph).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that reads exactly \"Commit pinkesh\" (ignore the\n// trailing paragraph-mark newline Office.js sometimes reports).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.replace(/[\\r\\n]+$/, \"\");\n  if (text === \"Commit pinkesh\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find a paragraph with text \"Commit pinkesh\"');\n}\n\n// Create the new (empty) paragraph directly after it \u2026\nconst newParagraph = target.insertParagraph(\"\", \"After\");\n\n// \u2026 then give it the same run / proofErr layout as the original\n// paragraph, with an extra \" 2\" run appended at the end:\n//   \"Commit \" + (spell-checked)\"pinkesh\" + \" 2\"\nconst newParagraphOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  '<w:r><w:t xml:space=\"preserve\">Commit </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>pinkesh</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> 2</w:t></w:r>' +\n  \"</w:p>\" +\n  \"</w:body></w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nnewParagraph.insertOoxml(newParagraphOoxml, \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that reads exactly \"Commit pinkesh\".\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text -replace \"[\\r\\n]+$\", \"\"\n    if ($t -eq \"Commit pinkesh\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw 'Could not find a paragraph with text \"Commit pinkesh\"'\n}\n\n# Insert a new (empty) paragraph right after it.\n$target.Range.InsertParagraphAfter()\n\n# Re-resolve that brand new paragraph via its index (one after $target).\n$newIndex = $target.Index + 1\n$newPara = $d.Paragraphs.Item($newIndex)\n\n# Give it the same run / proofErr layout as the original paragraph,\n# with an extra \" 2\" run appended at the end:\n#   \"Commit \" + (spell-checked)\"pinkesh\" + \" 2\"\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:t xml:space=\"preserve\">Commit </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>pinkesh</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> 2</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n\n$newPara.Range.InsertXML($xml)\n"}
